$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 128 ---
$ws.Range("L128").Value = 45912.658893564817
$ws.Range("M128").Value = 45912.658893356478
# Row 128 is no longer the last row; drop its special integer-style on E/F
$ws.Range("E128:F128").Style = "Normal"

# --- Append new rows 129-134 ---
# Row 129
$ws.Range("A129").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B129").Value = "MUTAG"
$ws.Range("C129").Value = "(5)-NN_Classifier_GED"
$ws.Range("D129").Value = 0.2
$ws.Range("E129").Value = 0.00000020557870370370369
$ws.Range("F129").Value = 0.00000006157407407407407
$ws.Range("G129").Value = 0.86842105263157898
$ws.Range("H129").Value = 0.8679574199770038
$ws.Range("I129").Value = 0.86872009569378005
$ws.Range("J129").Value = 0.86842105263157898
$ws.Range("K129").Value = 0.8641456582633052
$ws.Range("L129").Value = 45913.53986631944
$ws.Range("M129").Value = 45913.53986611111
$ws.Range("N129").Value = "Dummy_Calculator"
$ws.Range("O129").Value = "Simple Train-Test Split"
$ws.Range("L129:M129").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 130
$ws.Range("A130").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B130").Value = "MUTAG"
$ws.Range("C130").Value = "(5)-NN_Classifier_GED"
$ws.Range("D130").Value = 0.2
$ws.Range("E130").Value = 0.00000020557870370370369
$ws.Range("F130").Value = 0.00000006157407407407407
$ws.Range("G130").Value = 0.89473684210526316
$ws.Range("H130").Value = 0.89383971291866037
$ws.Range("I130").Value = 0.89809305873379097
$ws.Range("J130").Value = 0.89473684210526316
$ws.Range("K130").Value = 0.88795518207282908
$ws.Range("L130").Value = 45913.53986631944
$ws.Range("M130").Value = 45913.53986611111
$ws.Range("N130").Value = "Dummy_Calculator"
$ws.Range("O130").Value = "Hyperparameter Tuning (grid)"
$ws.Range("L130:M130").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 131
$ws.Range("A131").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B131").Value = "MUTAG"
$ws.Range("C131").Value = "(5)-NN_Classifier_GED"
$ws.Range("D131").Value = 0.2
$ws.Range("E131").Value = 0.0000002166550925925926
$ws.Range("F131").Value = 0.000000060578703703703707
$ws.Range("G131").Value = 0.86842105263157898
$ws.Range("H131").Value = 0.8635637494063296
$ws.Range("I131").Value = 0.87180451127819547
$ws.Range("J131").Value = 0.86842105263157898
$ws.Range("K131").Value = 0.8261538461538461
$ws.Range("L131").Value = 45913.617764120368
$ws.Range("M131").Value = 45913.617763900459
$ws.Range("N131").Value = "Dummy_Calculator"
$ws.Range("O131").Value = "Simple Train-Test Split"
$ws.Range("L131:M131").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 132
$ws.Range("A132").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B132").Value = "MUTAG"
$ws.Range("C132").Value = "(5)-NN_Classifier_GED"
$ws.Range("D132").Value = 0.2
$ws.Range("E132").Value = 0.0000002166550925925926
$ws.Range("F132").Value = 0.000000060578703703703707
$ws.Range("G132").Value = 0.8157894736842105
$ws.Range("H132").Value = 0.8173210804789752
$ws.Range("I132").Value = 0.8200187969924811
$ws.Range("J132").Value = 0.8157894736842105
$ws.Range("K132").Value = 0.8046153846153845
$ws.Range("L132").Value = 45913.617764120368
$ws.Range("M132").Value = 45913.617763900459
$ws.Range("N132").Value = "Dummy_Calculator"
$ws.Range("O132").Value = "Hyperparameter Tuning (grid)"
$ws.Range("L132:M132").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 133
$ws.Range("A133").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B133").Value = "MUTAG"
$ws.Range("C133").Value = "(5)-NN_Classifier_GED"
$ws.Range("D133").Value = 0.2
$ws.Range("E133").Value = 0.00000020427083333333329
$ws.Range("F133").Value = 0.00000006020833333333333
$ws.Range("G133").Value = 0.86842105263157898
$ws.Range("H133").Value = 0.8643962848297214
$ws.Range("I133").Value = 0.87293992557150446
$ws.Range("J133").Value = 0.86842105263157898
$ws.Range("K133").Value = 0.8363095238095238
$ws.Range("L133").Value = 45913.617978321759
$ws.Range("M133").Value = 45913.61797810185
$ws.Range("N133").Value = "GEDLIB_Calculator"
$ws.Range("O133").Value = "Simple Train-Test Split"
$ws.Range("L133:M133").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 134
$ws.Range("A134").Value = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$ws.Range("B134").Value = "MUTAG"
$ws.Range("C134").Value = "(5)-NN_Classifier_GED"
$ws.Range("D134").Value = 0.2
$ws.Range("E134").Value = 0.00000020427083333333329
$ws.Range("F134").Value = 0.00000006020833333333333
$ws.Range("G134").Value = 0.8157894736842105
$ws.Range("H134").Value = 0.79681020733652307
$ws.Range("I134").Value = 0.85738539898132426
$ws.Range("J134").Value = 0.8157894736842105
$ws.Range("K134").Value = 0.75
$ws.Range("L134").Value = 45913.617978318739
$ws.Range("M134").Value = 45913.61797810708
$ws.Range("N134").Value = "GEDLIB_Calculator"
$ws.Range("O134").Value = "Hyperparameter Tuning (grid)"
$ws.Range("L134:M134").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 134 (new last row) takes on the special integer number format for E/F
# that row 128 used to have before it was superseded
$ws.Range("E134:F134").NumberFormat = "0"
